# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels in AD1:AF1 ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting already used by the other header cells (bold,
# centered, bordered) by copying the style from the last existing header
# cell (AC1) onto the three new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows (rows 2-42): season record values for every player ---
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 74  # column AD = Wins
    $ws.Cells.Item($r, 31).Value = 88  # column AE = Losses
    $ws.Cells.Item($r, 32).Value = 0   # column AF = Ties
}

Write-Output "Added Wins/Losses/Ties columns for $($lastRow - 1) rows"
